$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing G14: "no" -> "confirm"
$ws.Range("G14").Value = "confirm"

# Copy the formatting of row 14 (A:G) down into the new row 15 so the new
# cells inherit the same cell styles (s="1" for A, s="2" for C/D, s="0" for
# B/E/F/G) before we overwrite the values.
$ws.Range("A14:G14").Copy()
$ws.Range("A15:G15").PasteSpecial(-4122)

# New row 15 values
$ws.Range("A15").Value = "com.hamxa.shaynachim"
$ws.Range("B15").Value = "bitcoin"
$ws.Range("C15").Value = "ronoren61@gmail.com"
$ws.Range("D15").Value = "nitanoren23@gmail.com"
$ws.Range("E15").Value = "27/5/2019 15:59"
$ws.Range("F15").Value = "amazing series of app. This and the blockchain app"
$ws.Range("G15").Value = "confirm"

# New hyperlinks for the email cells on row 15
$ws.Hyperlinks.Add($ws.Range("C15"), "mailto:ronoren61@gmail.com", "", "", "ronoren61@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D15"), "mailto:nitanoren23@gmail.com", "", "", "nitanoren23@gmail.com")

# Hyperlinks.Add swaps in Excel's built-in "Hyperlink" font style; re-apply
# the plain row style (same as C14/D14) so C15/D15 match the rest of the
# email columns, which don't carry a distinct hyperlink style either.
$ws.Range("C14:D14").Copy()
$ws.Range("C15:D15").PasteSpecial(-4122)

# Move the selection / active cell to A15 (matches the author's saved view state)
$null = $ws.Range("A15").Select()
